$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.532.44'
$ws.Range('E2').Value = '  +5.87%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.394.07'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '576.97'
$ws.Range('E5').Value = '  +7.59%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '155.39'
$ws.Range('E6').Value = '  +7.11%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.402.39'
$ws.Range('E8').Value = '  +6.56%  '
$ws.Range('E9').Value = '  -0.12%  '
$ws.Range('E10').Value = '  +2.28%  '
$ws.Range('E11').Value = '  +7.28%  '
$ws.Range('E12').Value = '  +0.92%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.981.48'
$ws.Range('E13').Value = '  +6.49%  '
$ws.Range('E14').Value = '  +0.40%  '
$ws.Range('E15').Value = '  +6.60%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '27.03'
$ws.Range('E16').Value = '  +4.80%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.605.91'
$ws.Range('E17').Value = '  +5.96%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.399.07'
$ws.Range('E18').Value = '  +6.32%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.36'
$ws.Range('E19').Value = '  +2.15%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.03'
$ws.Range('E20').Value = '  +6.07%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '8.42'
$ws.Range('E21').Value = '  +2.87%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '387.11'
$ws.Range('E22').Value = '  +4.89%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.998'
$ws.Range('E23').Value = '  -0.21%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.534'
$ws.Range('E24').Value = '  +2.15%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '70.80'
$ws.Range('E25').Value = '  +1.97%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.51'
$ws.Range('E26').Value = '  +10.92%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.181'
$ws.Range('E27').Value = '  +6.86%  '
$ws.Range('E28').Value = '  +17.89%  '
$ws.Range('E29').Value = '  +1.69%  '
$ws.Range('E30').Value = '  +8.14%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.57'
$ws.Range('E31').Value = '  +7.55%  '
$ws.Range('E32').Value = '  +13.57%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.63'
$ws.Range('E33').Value = '  +6.92%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '23.12'
$ws.Range('E34').Value = '  +2.89%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.70'
$ws.Range('E35').Value = '  +1.91%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.48'
$ws.Range('E36').Value = '  +9.49%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '158.49'
$ws.Range('E37').Value = '  +0.60%  '
$ws.Range('B38').Value = 'Stacks'
$ws.Range('C38').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.87'
$ws.Range('E38').Value = '  +10.83%  '
$ws.Range('B39').Value = 'EnergySwap'
$ws.Range('C39').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '27.44'
$ws.Range('E39').Value = '  +3.80%  '
$ws.Range('E40').Value = '  +7.67%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.898.11'
$ws.Range('E41').Value = '  +3.94%  '
$ws.Range('E42').Value = '  +4.33%  '
$ws.Range('E43').Value = '  +6.22%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '41.44'
$ws.Range('E44').Value = '  +4.11%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.30'
$ws.Range('E45').Value = '  +1.98%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.06'
$ws.Range('E46').Value = '  +8.29%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.443.57'
$ws.Range('E47').Value = '  +6.64%  '
$ws.Range('E48').Value = '  +8.46%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '299.52'
$ws.Range('E49').Value = '  +14.27%  '
$ws.Range('E50').Value = '  -2.12%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.31'
$ws.Range('E51').Value = '  +2.72%  '
